$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the header info (student name / id) that used to live on row 1.
$ws.Range("C1:E1").ClearContents()

# --- Fall/Spring/Summer 2022 block (rows 4-9) ---
# Fall 2022 (col A/B) and Spring 2022 (col C/D) get new course lists; the
# Summer 2022 column (E/F) is removed entirely for this block.
$ws.Range("A4").Value = "POLS 1101"
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = "CPSC 3165"
$ws.Range("D4").Value = 3
$ws.Range("E4:F4").ClearContents()

$ws.Range("A5").Value = "DSCI 3111"
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = "CPSC 3415"
$ws.Range("D5").Value = 1
$ws.Range("E5:F5").ClearContents()

$ws.Range("A6").Value = "ARTH 3117"
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = "CPSC 4135"
$ws.Range("D6").Value = 3

$ws.Range("A7").Value = "CPSC 3121"
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = "CPSC 4148"
$ws.Range("D7").Value = 3

$ws.Range("A8").Value = "ARTH 3136"
$ws.Range("B8").Value = 3
$ws.Range("C8").Value = "CPSC 4155"
$ws.Range("D8").Value = 3

$ws.Range("A9").Value = "CPSC 4000"
$ws.Range("B9").Value = 0

# --- Fall/Spring/Summer 2023 block (rows 13-15) ---
# Fall 2023 (col A/B) gets a new course list; Spring 2023 (col C/D) shrinks
# to a single row; Summer 2023 column (E/F) is removed entirely.
$ws.Range("A13").Value = "CPSC 4157"
$ws.Range("B13").Value = 3
$ws.Range("C13").Value = "CPSC 4176"
$ws.Range("D13").Value = 3
$ws.Range("E13:F13").ClearContents()

$ws.Range("A14").Value = "CPSC 4175"
$ws.Range("B14").Value = 3
$ws.Range("C14:D14").ClearContents()

$ws.Range("A15").Value = "CPSC 4555"
$ws.Range("B15").Value = 3
$ws.Range("C15:D15").ClearContents()
